# Apply the checklist-grading updates to the first sheet
# ("Rúbrica con Evidencias-1ra") of the workbook: a couple of indicator
# levels move from 2 -> 3, and the accompanying reviewer-comment notes in
# column I are revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7 (Objetivo / Producto indicator) ---------------------------------
# Level raised from 2 to 3, and the comment is replaced.
$ws.Cells.Item(7, 8).Value = 3
$ws.Cells.Item(7, 9).Value = "Alguna limitación?"

# --- Row 11 (Requisitos Funcionales indicator) ------------------------------
# Level raised from 2 to 3; the old comment no longer applies, so it is
# cleared out entirely.
$ws.Cells.Item(11, 8).Value = 3
$ws.Cells.Item(11, 9).Value = ""

# --- Row 13 (Priorización indicator) ----------------------------------------
# Comment text revised (level stays at 2). The note now uses a smaller,
# vertically-centered, wrapped font to fit the cell better.
$i13 = $ws.Cells.Item(13, 9)
$i13.Value = "En la tabla de gestión no se tiene el dato de la prioridad. Lo consideraron en algún otro documento?"
$i13.Font.Size = 6
$i13.Font.Name = "Arial"
$i13.WrapText = $true
$i13.VerticalAlignment = -4108

# --- Row 14 (Artefactos indicator) ------------------------------------------
# Comment text revised (level stays at 3). Reuses the exact same restyle as
# row 13 (same font/alignment), applied via copy/paste-special of the
# formatting so both cells end up sharing one single cell style.
$i14 = $ws.Cells.Item(14, 9)
$i14.Value = "Criterios de aceptación"
$i13.Copy()
$i14.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Keep the saved cursor/selection in sync with the edited area ----------
$ws.Range("I14").Select()
